$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (Volume/Number and report date range)
$ws.Range("A8").Value = "Volume 29   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/28/2022  Through  12/4/2022"

# Crime-complaint grid updates
$ws.Range("G14").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("G14").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Value = "'***.*"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null
$ws.Range("F15").Value = 4
$ws.Range("M15").Value = 143.75
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 30
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 172.727272727273
$ws.Range("I16").Value = 288
$ws.Range("J16").Value = 135
$ws.Range("K16").Value = 113.333333333333
$ws.Range("L16").Value = 60.893854748603
$ws.Range("M16").Value = 47.692307692307
$ws.Range("N16").Value = -74.331550802139
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 115.384615384615
$ws.Range("I17").Value = 317
$ws.Range("J17").Value = 191
$ws.Range("K17").Value = 65.968586387434
$ws.Range("L17").Value = 26.8
$ws.Range("M17").Value = 54.634146341463
$ws.Range("N17").Value = -2.16049382716
$ws.Range("C18").Value = 14
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = 75
$ws.Range("F18").Value = 62
$ws.Range("G18").Value = 50
$ws.Range("H18").Value = 24
$ws.Range("I18").Value = 540
$ws.Range("J18").Value = 339
$ws.Range("K18").Value = 59.29203539823
$ws.Range("L18").Value = 33.995037220843
$ws.Range("M18").Value = 20.267260579064
$ws.Range("N18").Value = -75.838926174496
$ws.Range("C19").Value = 33
$ws.Range("D19").Value = 162
$ws.Range("E19").Value = -79.629629629629
$ws.Range("F19").Value = 117
$ws.Range("G19").Value = 216
$ws.Range("H19").Value = -45.833333333333
$ws.Range("I19").Value = 1400
$ws.Range("J19").Value = 828
$ws.Range("K19").Value = 69.082125603864
$ws.Range("L19").Value = 66.865315852205
$ws.Range("M19").Value = 124.719101123596
$ws.Range("N19").Value = 10.062893081761
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 31
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 72.222222222222
$ws.Range("I20").Value = 260
$ws.Range("J20").Value = 188
$ws.Range("K20").Value = 38.297872340425
$ws.Range("L20").Value = 75.675675675675
$ws.Range("M20").Value = -5.109489051094
$ws.Range("N20").Value = -93.252011419673
$ws.Range("C21").Value = 65
$ws.Range("D21").Value = 181
$ws.Range("E21").Value = -64.088397790055
$ws.Range("F21").Value = 274
$ws.Range("G21").Value = 308
$ws.Range("H21").Value = -11.038961038961
$ws.Range("I21").Value = 2850
$ws.Range("J21").Value = 1703
$ws.Range("K21").Value = 67.351732237228
$ws.Range("L21").Value = 53.804641122504
$ws.Range("M21").Value = 61.107970604861
$ws.Range("N21").Value = -67.78204838345
$ws.Range("D22").Value = 1
$ws.Range("C16").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = 0
$ws.Range("H16").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 38
$ws.Range("J22").Value = 14
$ws.Range("K22").Value = 171.428571428571
$ws.Range("L22").Value = 1166.66666666667
$ws.Range("M22").Value = 850
$ws.Range("M23").Value = 45.454545454545
$ws.Range("C24").Value = 53
$ws.Range("D24").Value = 67
$ws.Range("E24").Value = -20.895522388059
$ws.Range("F24").Value = 232
$ws.Range("G24").Value = 221
$ws.Range("H24").Value = 4.97737556561
$ws.Range("I24").Value = 2590
$ws.Range("J24").Value = 1985
$ws.Range("K24").Value = 30.478589420654
$ws.Range("L24").Value = 55.368926214757
$ws.Range("M24").Value = 75.831636116768
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 53
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 39.473684210526
$ws.Range("I25").Value = 660
$ws.Range("J25").Value = 487
$ws.Range("K25").Value = 35.523613963039
$ws.Range("L25").Value = 34.419551934826
$ws.Range("M25").Value = -1.639344262295
$ws.Range("D26").Value = 1
$ws.Range("C16").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Value = -100
$ws.Range("H16").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 150
$ws.Range("J26").Value = 29
$ws.Range("K26").Value = 75.862068965517
$ws.Range("D27").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value = "'***.*"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 8
$ws.Range("L27").Value = 34.482758620689
$ws.Range("C28").Value = 1
$ws.Range("C16").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 8
$ws.Range("K28").Value = 14.285714285714
$ws.Range("L28").Value = 14.285714285714
$ws.Range("M28").Value = 166.666666666667
$ws.Range("N28").Value = -68
$ws.Range("C29").Value = 1
$ws.Range("C16").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null
$ws.Range("F29").Value = 3
$ws.Range("H29").Value = 200
$ws.Range("I29").Value = 8
$ws.Range("K29").Value = 33.333333333333
$ws.Range("L29").Value = 14.285714285714
$ws.Range("M29").Value = 166.666666666667
$ws.Range("N29").Value = -61.904761904761
